$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("C5").Value = 212.3897266062013
$ws.Range("D5").Value = 45109.39596785694
$ws.Range("E5").Value = 117.1594356773392
$ws.Range("F5").Value = 0.7502542900481672
$ws.Range("G5").Value = 66.71504095686809
$ws.Range("H5").Value = 34.06838452029956
$ws.Range("I5").Value = 44.58809855828566

# Row 6
$ws.Range("C6").Value = 252.9315743791161
$ws.Range("D6").Value = 63974.38131789833
$ws.Range("E6").Value = 144.7814780756655
$ws.Range("F6").Value = 0.64727315429288
$ws.Range("G6").Value = 117.9797277160395
$ws.Range("H6").Value = 44.39718861007679
$ws.Range("I6").Value = 58.15845229744227

# Row 7
$ws.Range("C7").Value = 223.2176590451196
$ws.Range("D7").Value = 49826.12330958326
$ws.Range("E7").Value = 124.0656505615306
$ws.Range("F7").Value = 0.7244263828166845
$ws.Range("G7").Value = 79.53251975239168
$ws.Range("H7").Value = 35.97821434897747
$ws.Range("I7").Value = 47.9810329990345

# Row 8
$ws.Range("C8").Value = 213.2665011982465
$ws.Range("D8").Value = 45482.60053334165
$ws.Range("E8").Value = 117.5513989429431
$ws.Range("F8").Value = 0.7481880633305527
$ws.Range("G8").Value = 65.18961038028252
$ws.Range("H8").Value = 33.52023773910638
$ws.Range("I8").Value = 44.10296707848546

# Row 9
$ws.Range("C9").Value = 252.3787475920975
$ws.Range("D9").Value = 63695.03223615564
$ws.Range("E9").Value = 144.3897344274109
$ws.Range("F9").Value = 0.6488133633331933
$ws.Range("G9").Value = 112.8983413844712
$ws.Range("H9").Value = 44.13286374636161
$ws.Range("I9").Value = 57.54215248169891

# Row 10
$ws.Range("C10").Value = 223.687668021671
$ws.Range("D10").Value = 50036.17282497331
$ws.Range("E10").Value = 124.2616671163417
$ws.Range("F10").Value = 0.7232646607942026
$ws.Range("G10").Value = 77.11800957016203
$ws.Range("H10").Value = 35.8989187199584
$ws.Range("I10").Value = 47.46310609082383

# Row 11
$ws.Range("C11").Value = 211.8131792658207
$ws.Range("D11").Value = 44864.82291069469
$ws.Range("E11").Value = 116.6528725909557
$ws.Range("F11").Value = 0.7516083554371076
$ws.Range("G11").Value = 64.16376049440402
$ws.Range("H11").Value = 33.7394191518917
$ws.Range("I11").Value = 43.90355349083377

# Row 12
$ws.Range("C12").Value = 252.9646531078764
$ws.Range("D12").Value = 63991.11572198826
$ws.Range("E12").Value = 144.9161559779944
$ws.Range("F12").Value = 0.6471808880849412
$ws.Range("G12").Value = 111.4284350645573
$ws.Range("H12").Value = 44.00129555714123
$ws.Range("I12").Value = 57.53724510622013

# Row 13
$ws.Range("C13").Value = 222.8158068457923
$ws.Range("D13").Value = 49646.88378034143
$ws.Range("E13").Value = 123.719414072172
$ws.Range("F13").Value = 0.7254177038774965
$ws.Range("G13").Value = 75.98113425358103
$ws.Range("H13").Value = 35.85926827662228
$ws.Range("I13").Value = 47.31232401557826

# Row 14
$ws.Range("C14").Value = 212.5515720898156
$ws.Range("D14").Value = 45178.17079785207
$ws.Range("E14").Value = 116.9130617880801
$ws.Range("F14").Value = 0.7498735219537288
$ws.Range("G14").Value = 64.17614742621032
$ws.Range("H14").Value = 33.70204931285453
$ws.Range("I14").Value = 44.01034708315365

# Row 15
$ws.Range("C15").Value = 253.2008638651642
$ws.Range("D15").Value = 64110.6774620654
$ws.Range("E15").Value = 144.8372007124044
$ws.Range("F15").Value = 0.6465216767791671
$ws.Range("G15").Value = 112.1460990039648
$ws.Range("H15").Value = 44.10102726131703
$ws.Range("I15").Value = 57.54461255958704

# Row 16
$ws.Range("C16").Value = 223.4094451667733
$ws.Range("D16").Value = 49911.7801897255
$ws.Range("E16").Value = 123.8948085063851
$ws.Range("F16").Value = 0.7239526398334952
$ws.Range("G16").Value = 76.16985841987325
$ws.Range("H16").Value = 35.88469233783497
$ws.Range("I16").Value = 47.39425853807221

# Row 17
$ws.Range("C17").Value = 212.3280355796343
$ws.Range("D17").Value = 45083.19469310644
$ws.Range("E17").Value = 116.8144619027489
$ws.Range("F17").Value = 0.7503993519764819
$ws.Range("G17").Value = 64.17954578787513
$ws.Range("H17").Value = 33.6373967881255
$ws.Range("I17").Value = 44.00948267542442

# Row 18
$ws.Range("C18").Value = 253.007474305661
$ws.Range("D18").Value = 64012.78205452974
$ws.Range("E18").Value = 144.7580938057344
$ws.Range("F18").Value = 0.6470614293738453
$ws.Range("G18").Value = 111.847245045968
$ws.Range("H18").Value = 43.97466872776717
$ws.Range("I18").Value = 57.54480265205322

# Row 19
$ws.Range("C19").Value = 223.1951482108165
$ws.Range("D19").Value = 49816.07418484833
$ws.Range("E19").Value = 123.8010823627355
$ws.Range("F19").Value = 0.7244819615667214
$ws.Range("G19").Value = 76.09768599503626
$ws.Range("H19").Value = 35.78464321974565
$ws.Range("I19").Value = 47.39365778227862

